# Add a new "2020" column (Q) to the worksheet, mirroring the existing
# column P (2019) formatting/styles, and fill in the new year's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: thin separator row - Q3 just needs P3's (empty) style, no value.
$ws.Range("P3").Copy($ws.Range("Q3"))

# Row 4 (year headers): Q4 = 2020, same style as P4.
$ws.Range("P4").Copy($ws.Range("Q4"))
$ws.Range("Q4").Value = 2020

# Row 5: new data value, same style as P5.
$ws.Range("P5").Copy($ws.Range("Q5"))
$ws.Range("Q5").Value = 25.6

# Row 6: new data value, same style as P6.
$ws.Range("P6").Copy($ws.Range("Q6"))
$ws.Range("Q6").Value = 13.073527219449954

# Row 7: new data value, same style as P7.
$ws.Range("P7").Copy($ws.Range("Q7"))
$ws.Range("Q7").Value = 21.941290626870046

# Row 8: new data value, same style as P8.
$ws.Range("P8").Copy($ws.Range("Q8"))
$ws.Range("Q8").Value = 196.6

# Row 1 picks up a custom height in the source diff.
$ws.Rows.Item(1).RowHeight = 19.5
